# feat: add 2022-Q1 data
#
# Adds a new "2022-Q1" worksheet (per-fund holding detail, same shape as the
# existing quarter sheets) positioned right after "2021-Q4", and updates the
# "总计" (grand-total) summary sheet with a new leading row for 2022-Q1
# (9 funds held, 1.98 亿元 total market value), pushing the older rows down.

$wb = $excel.ActiveWorkbook

function Set-TextCell($range, [string]$value) {
    # Force text storage so numeric-looking strings (fund codes with leading
    # zeros, "24.64", "0.8452", ...) are never silently coerced to numbers.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Sheet whose cell formatting (borders / bold / alignment) we clone for the
# new sheet's header row + index column, so we don't hand-roll new styles.
$fmtSrc = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------------
# 1. Remove the existing "总计" sheet - it gets rebuilt (with the new row)
#    below so it ends up last again, after the freshly inserted "2022-Q1".
# ---------------------------------------------------------------------------
$oldTotal = $wb.Worksheets.Item("总计")
$oldTotal.Delete()

# ---------------------------------------------------------------------------
# 2. New "2022-Q1" sheet (fund-level holdings), inserted right after
#    "2021-Q4".
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $q4)
$q1.Name = "2022-Q1"

Set-TextCell $q1.Range("B1") "基金代码"
Set-TextCell $q1.Range("C1") "基金名称"
Set-TextCell $q1.Range("D1") "基金规模"
Set-TextCell $q1.Range("E1") "股票总仓位"
Set-TextCell $q1.Range("F1") "仓位占比"
Set-TextCell $q1.Range("G1") "持有市值(亿元)"
Set-TextCell $q1.Range("H1") "仓位排名"
$fmtSrc.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$q1Rows = @(
    @{code="920003"; name="中金新锐股票A";                   scale="24.64"; pos="92.76"; ratio="3.43"; mv="0.8452"; rank=6},
    @{code="163804"; name="中银收益混合A";                   scale="19.19"; pos="85.56"; ratio="3.21"; mv="0.6160"; rank=4},
    @{code="163822"; name="中银主题策略混合";                 scale="7.89";  pos="84.73"; ratio="2.93"; mv="0.2312"; rank=8},
    @{code="920923"; name="中金新锐股票C";                   scale="3.94";  pos="92.76"; ratio="3.43"; mv="0.1351"; rank=6},
    @{code="200010"; name="长城双动力混合";                   scale="2.03";  pos="86.77"; ratio="4.69"; mv="0.0952"; rank=2},
    @{code="014505"; name="中银收益混合C";                   scale="0.98";  pos="85.56"; ratio="3.21"; mv="0.0315"; rank=4},
    @{code="002292"; name="诺安益鑫灵活配置混合";              scale="0.30";  pos="50.08"; ratio="4.46"; mv="0.0134"; rank=4},
    @{code="004320"; name="前海开源沪港深乐享生活灵活配置混合"; scale="0.17";  pos="85.24"; ratio="5.87"; mv="0.0100"; rank=3},
    @{code="960012"; name="中银收益混合H";                   scale="0.03";  pos="85.56"; ratio="3.21"; mv="0.0010"; rank=4}
)

for ($i = 0; $i -lt $q1Rows.Count; $i++) {
    $r = $i + 2
    $row = $q1Rows[$i]
    $q1.Range("A$r").Value = $i
    Set-TextCell $q1.Range("B$r") $row.code
    Set-TextCell $q1.Range("C$r") $row.name
    Set-TextCell $q1.Range("D$r") $row.scale
    Set-TextCell $q1.Range("E$r") $row.pos
    Set-TextCell $q1.Range("F$r") $row.ratio
    Set-TextCell $q1.Range("G$r") $row.mv
    $q1.Range("H$r").Value = $row.rank
}

$fmtSrc.Range("A2").Copy()
$q1.Range("A2:A10").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Rebuilt "总计" sheet, inserted right after "2022-Q1" (i.e. last again),
#    with the new 2022-Q1 totals row on top and the previous rows shifted
#    down by one.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

Set-TextCell $total.Range("B1") "日期"
Set-TextCell $total.Range("C1") "持有数量(只)"
Set-TextCell $total.Range("D1") "持有市值(亿元)"
$fmtSrc.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$totalRows = @(
    @{label="2022-Q1"; count=9; mv=1.98},
    @{label="2021-Q4"; count=7; mv=5.06},
    @{label="2021-Q2"; count=4; mv=1.11},
    @{label="2021-Q1"; count=2; mv=0.02},
    @{label="2020-Q4"; count=2; mv=2.57}
)

for ($i = 0; $i -lt $totalRows.Count; $i++) {
    $r = $i + 2
    $row = $totalRows[$i]
    $total.Range("A$r").Value = $i
    Set-TextCell $total.Range("B$r") $row.label
    $total.Range("C$r").Value = $row.count
    $total.Range("D$r").Value = $row.mv
}

$fmtSrc.Range("A2").Copy()
$total.Range("A2:A6").PasteSpecial(-4122)
